# Auto-update draw results: append the 2025-10-22 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data lives in rows 1..35 (row 1 = header); the new draw goes
# immediately below the last populated row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$target = $ws.Range("A" + $newRow + ":E" + $newRow)

# The existing rows store every value (dates, phase codes, results,
# timestamps) as literal text, not as numbers/dates. Force this new row to
# Text format first so Excel doesn't auto-coerce the numeric-looking values
# ("2025-10-22", "251022") into a date serial / number.
$target.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-10-22"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "251022"
$ws.Cells.Item($newRow, 4).Value = "4-7-2-6"
$ws.Cells.Item($newRow, 5).Value = "2025-10-22T21:37:36.319+04:00"

# Re-apply the formatting of the previous data row (which uses the
# workbook's plain default style) onto the new row so we don't leave the
# "@" text-format override (and its extra style record) behind.
$ws.Range("A" + $lastRow + ":E" + $lastRow).Copy()
$target.PasteSpecial(-4122)
$excel.CutCopyMode = $false
